$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the header label currently stored in W5 (shared string "nomor_s")
$headerLabel = $ws.Range("W5").Value2

# Remove the old column W cells (W5:W119) completely, content and formatting,
# so the <c> nodes disappear entirely from the XML.
$ws.Range("W5:W119").Clear()

# Re-create the same data one column to the right (column X), with the
# numeric value changed from 2070 to 2041.
$ws.Range("X5").Value = $headerLabel
$ws.Range("X6").Value = 2041
for ($r = 7; $r -le 119; $r++) {
    $ws.Cells.Item($r, 24).Value = 2041
}
